$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 108, pushing the existing rows 108-128 down to 109-129.
$ws.Rows.Item(108).Insert()

# Populate the newly inserted row 108 with the new daily price record (same fixed
# columns as the rest of the "Coco" block, with this week's volume/price figures).
$ws.Range("A108").Value = 10
$ws.Range("B108").Value = "Vega Modelo de Temuco"
$ws.Range("C108").Value = "La Araucanía"
$ws.Range("D108").Value = 45244
$ws.Range("D108").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E108").Value = 9
$ws.Range("F108").Value = "Fruta"
$ws.Range("G108").Value = 100108
$ws.Range("H108").Value = "Tropicales y subtropicales"
$ws.Range("I108").Value = 100108007
$ws.Range("J108").Value = "Coco"
$ws.Range("K108").Value = "Sin especificar"
$ws.Range("L108").Value = "Primera"
$ws.Range("M108").Value = 50
$ws.Range("N108").Value = 35000
$ws.Range("O108").Value = 35000
$ws.Range("P108").Value = 35000
$ws.Range("Q108").Value = "$/malla 20 unidades"
$ws.Range("R108").Value = "Perú"
$ws.Range("S108").Value = 1750
$ws.Range("T108").Value = 20
